$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (column D) and volume-change (column E) values.
# Column D values are forced to Text format before assignment so that Excel does not
# auto-convert numeric-looking strings (e.g. "27.624.42") into floating point numbers,
# then the style is reset back to Normal so no residual formatting is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.624.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.667.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.26%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.508'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.17'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.264'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.904.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.676.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.76%  '
$ws.Range("E14").Value = '  -3.11%  '
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.616.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("E19").Value = '  -2.95%  '
$ws.Range("E20").Value = '  -3.96%  '
$ws.Range("E22").Value = '  -2.41%  '
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("E24").Value = '  -3.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("E26").Value = '  -3.66%  '
$ws.Range("E27").Value = '  -1.19%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("E30").Value = '  +3.39%  '
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.468.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("E34").Value = '  -4.46%  '
$ws.Range("E35").Value = '  -5.09%  '
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("E37").Value = '  -3.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.579'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.76%  '
$ws.Range("E39").Value = '  -1.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '69.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.49%  '
$ws.Range("E41").Value = '  -4.43%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.51%  '
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.811.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.790'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.18%  '
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("E51").Value = '  -4.04%  '
